# Append 13 new game rows (897-909) to Sheet1, matching new box-score data
# pulled into the "ballgorithm" NBA 2023-24 tracker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Atlanta Hawks',102,'Brooklyn Nets',114,'No',17832,'Barclays Center','Brooklyn Nets','Atlanta Hawks'),
    @('Utah Jazz',120,'Miami Heat',126,'No',17832,'Kaseya Center','Miami Heat','Utah Jazz'),
    @('Portland Trail Blazers',107,'Memphis Grizzlies',100,'OT',17832,'FedEx Forum','Portland Trail Blazers','Memphis Grizzlies'),
    @('Denver Nuggets',124,'Los Angeles Lakers',114,'No',17832,'Crypto.com Arena','Denver Nuggets','Los Angeles Lakers'),
    @('Houston Rockets',118,'Phoenix Suns',109,'No',17832,'Footprint Center','Houston Rockets','Phoenix Suns'),
    @('Philadelphia 76ers',120,'Dallas Mavericks',116,'No',17832,'American Airlines Center','Philadelphia 76ers','Dallas Mavericks'),
    @('Golden State Warriors',88,'Boston Celtics',140,'No',17832,'TD Garden','Boston Celtics','Golden State Warriors'),
    @('Los Angeles Clippers',89,'Minnesota Timberwolves',88,'No',17832,'Target Center','Los Angeles Clippers','Minnesota Timberwolves'),
    @('Detroit Pistons',91,'Orlando Magic',113,'No',17832,'Amway Center','Orlando Magic','Detroit Pistons'),
    @('Charlotte Hornets',106,'Toronto Raptors',111,'No',17832,'Scotiabank Arena','Toronto Raptors','Charlotte Hornets'),
    @('New York Knicks',107,'Cleveland Cavaliers',98,'No',17832,'Rocket Mortgage Fieldhouse','New York Knicks','Cleveland Cavaliers'),
    @('Indiana Pacers',105,'San Antonio Spurs',117,'No',17832,'Frost Bank Center','San Antonio Spurs','Indiana Pacers'),
    @('Oklahoma City Thunder',118,'Phoenix Suns',110,'No',17832,'Footprint Center','Oklahoma City Thunder','Phoenix Suns')
)

$startRow = 897
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Range("A$r").Value = $rowVals[0]
    $ws.Range("B$r").Value = $rowVals[1]
    $ws.Range("C$r").Value = $rowVals[2]
    $ws.Range("D$r").Value = $rowVals[3]
    $ws.Range("E$r").Value = $rowVals[4]
    $ws.Range("F$r").Value = $rowVals[5]
    $ws.Range("G$r").Value = $rowVals[6]
    $ws.Range("H$r").Value = $rowVals[7]
    $ws.Range("I$r").Value = $rowVals[8]
}

$lastRow = $startRow + $data.Count - 1

# Scroll the view down and select G897, mirroring the author's on-screen state
$excel.ActiveWindow.ScrollRow = 877
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G897").Select() | Out-Null
